$wb = $excel.ActiveWorkbook

# --- Sheet "Info": update selection only (no data changes) ---
$wsInfo = $wb.Worksheets.Item("Info")
$wsInfo.Range("D12").Select()

# --- Sheet "A1": fill evidence row ---
$wsA1 = $wb.Worksheets.Item("A1")
$wsA1.Range("A2").Value = "79DAEDA84C1D760394B3B579B70B0379B94368AE9006E663C249D733E65AC8F0"
$wsA1.Range("B2").Value = "smtnftcoll"
$wsA1.Range("B2").Select()

# --- Sheet "A2": fill evidence rows (row2 + new row3) ---
$wsA2 = $wb.Worksheets.Item("A2")
$wsA2.Range("A2").Value = "1B557A1D5C5956E4B31C3B2F075D95BBB2F38A3978FDB13095BC3C29440D1855"
$wsA2.Range("B2").Value = "smtnftcoll"
$wsA2.Range("C2").Value = "nft0001"
$wsA2.Range("A3").Value = "B70D2752DD3489652D059C5048D36AD767835179B16C3B52C63FCF9E15D562B6"
$wsA2.Range("B3").Value = "smtnftcoll"
$wsA2.Range("C3").Value = "nft0002"
$wsA2.Range("C3").Select()

# --- Sheet "A3": fill evidence row ---
$wsA3 = $wb.Worksheets.Item("A3")
$wsA3.Range("A2").Value = "25C4BC866A388160FCDA5F6ED2C074FEE91974E88F7E4C869526DC8FEC104642"
$wsA3.Range("B2").Value = "juno1qm2k2nc0c56hyhxuc8kamwqaw75se7hw7whqqaqpvsmuvxjxk4usrnl053"
$wsA3.Range("C2").Value = "nft0001"
$wsA3.Range("D2").Value = "uni-6"
$wsA3.Range("B2:D2").Select()

# --- Sheet "A4": fill evidence row ---
$wsA4 = $wb.Worksheets.Item("A4")
$wsA4.Range("A2").Value = "1CE694CCEADD3EA4EF2FBBB9D6D1FB61A64544E57DB0D4E68FFA402020EFBAC7"
$wsA4.Range("B2").Value = "ibc/448DC656EA7119B20AEEEE9FA9D68182BA51C07FECD2A3BC0048C9FECD6D056D"
$wsA4.Range("C2").Value = "nft0002"
$wsA4.Range("D2").Value = "gon-flixnet-1"
$wsA4.Range("B2:D2").Select()

# --- Sheet "A5": fill evidence row ---
$wsA5 = $wb.Worksheets.Item("A5")
$wsA5.Range("A2").Value = "92908FFCE8A62F8059AB91790C3F9B3CCC9D12A592B60E5B4230FAD1B3DFBF7F"
$wsA5.Range("B2").Value = "juno1qm2k2nc0c56hyhxuc8kamwqaw75se7hw7whqqaqpvsmuvxjxk4usrnl053"
$wsA5.Range("C2").Value = "nft0001"
$wsA5.Range("D2").Value = "uni-6"
$wsA5.Range("B2:D2").Select()

# --- Sheet "A6": fill evidence row and make it the active tab ---
$wsA6 = $wb.Worksheets.Item("A6")
$wsA6.Range("A2").Value = "3AFAA48CD81221023882C045A24621C5212600B591A944BA2B8CD893AA6254A2"
$wsA6.Range("B2").Value = "ibc/448DC656EA7119B20AEEEE9FA9D68182BA51C07FECD2A3BC0048C9FECD6D056D"
$wsA6.Range("C2").Value = "nft0002"
$wsA6.Range("D2").Value = "gon-flixnet-1"
$wsA6.Activate()
$wsA6.Range("B2:D2").Select()
